$d = $word.ActiveDocument

# Locate the end of the "Teste" paragraph (the last paragraph in the body)
# and collapse a range to right after it.
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd = 0

# Insert a new paragraph whose text run carries no direct character
# formatting, but whose paragraph mark (pPr/rPr) is underlined - this
# mirrors selecting just the pilcrow and toggling Underline in the UI.
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
          '<w:p>' + `
            '<w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' + `
            '<w:r><w:t>Instalada tds os softwares NOW!!!!</w:t></w:r>' + `
          '</w:p>' + `
        '</w:body>' + `
      '</w:document>' + `
    '</pkg:xmlData></pkg:part>' + `
  '</pkg:package>'

$null = $endRange.InsertXML($xmlFrag)
